# Fortschritt.xlsx update ("Liederupdate + Fortschritt updated")
#
# Marks four additional songs (rows 9, 11, 24, 82 on Tabelle1) as part of
# the "Hymnen" (anthems) category and flags their progress as "ok" in
# columns B/C, mirroring the existing rows that already carry that
# annotation (e.g. row 4). Also updates the active selection to reflect
# where the editor left off.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 9 - "Brüder, reicht die Hand zum Bunde – Art Bundeshymne"
$ws.Range("B9").Value = "Hymnen"
$ws.Range("C9").Value = "ok"

# Row 11 - "Bundeshymne"
$ws.Range("B11").Value = "Hymnen"
$ws.Range("C11").Value = "ok"

# Row 24 - "Dort, wo des Wienerwaldes liebes Rauschen"
$ws.Range("B24").Value = "Hymnen"
$ws.Range("C24").Value = "ok"

# Row 82 - "Wien, mein Wien, gar oft besungen"
$ws.Range("B82").Value = "Hymnen"
$ws.Range("C82").Value = "ok"

# Reflect the editor's final cursor position/selection on the sheet.
[void]$ws.Range("B83").Select()
